# Updated diagram
#
# 1) Refresh the cached "datetimeFigureOut" date field text (shown on the
#    slide master and on every slide layout's Date Placeholder) from
#    12/21/2021 to 6/27/2022.
# 2) Rename the "Effective Activations" label box on slide 1 to
#    "Effective Notifications".

$p = $ppt.ActivePresentation

$oldDate = "12/21/2021"
$newDate = "6/27/2022"

# --- Slide master's own Date Placeholder -------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout's Date Placeholder ------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 1: "Effective Activations" -> "Effective Notifications" ------
$slide1 = $p.Slides.Item(1)
$topShape = $slide1.Shapes.Item(1)
$items = $topShape.GroupItems
for ($gi = 1; $gi -le $items.Count; $gi++) {
    $item = $items.Item($gi)
    if ($item.Name -eq "Rectangle 7") {
        if ($item.TextFrame.TextRange.Text -eq "Effective Activations") {
            $item.TextFrame.TextRange.Text = "Effective Notifications"
        }
    }
}
